# Atualização de bases das ligas, do dia: 08-04-2024 às 21:28
# Adds 3 new match rows (146-148) to the "Croatia HNL" sheet, mirroring the
# layout/formatting of the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{
        row = 146; id = 144; matchId = 6788942; date = 45395.49305555555
        home = "NK Varazdin"; away = "NK Lokomotiva Zagreb"
        K = 2.8;  L = 3.25; M = 2.5;   N = 2.8;  O = 3.25; P = 2.5
        Q = 0;    R = 2.05; S = 1.8;   T = 2.5;  U = 1.975; V = 1.875
        W = 0;    X = 0;    Y = 0;     Z = 0;    AA = 0
    },
    @{
        row = 147; id = 145; matchId = 6788943; date = 45395.58333333334
        home = "HNK Gorica"; away = "Dinamo Zagreb"
        K = 8;    L = 4.5;  M = 1.363; N = 8;    O = 4.5;  P = 1.363
        Q = 1.25; R = 2.025; S = 1.825; T = 2.75; U = 1.925; V = 1.925
        W = 0;    X = 0;    Y = 0;     Z = 0;    AA = 0
    },
    @{
        row = 148; id = 146; matchId = 6788944; date = 45396.58333333334
        home = "Istra 1961"; away = "HNK Rijeka"
        K = 5.5;  L = 3.6;  M = 1.615; N = 5.5;  O = 3.6;  P = 1.615
        Q = 0.75; R = 2.025; S = 1.825; T = 2.25; U = 1.925; V = 1.925
        W = 0;    X = 0;    Y = 0;     Z = 0;    AA = 0
    }
)

foreach ($r in $rows) {
    $n = $r.row

    $ws.Cells.Item($n, 1).Value = $r.id          # A - id
    $ws.Cells.Item($n, 2).Value = $r.matchId      # B
    $ws.Cells.Item($n, 3).Value = "Croatia HNL"   # C - Div
    $ws.Cells.Item($n, 4).Value = "Croatia HNL"   # D - Div Original Name
    $ws.Cells.Item($n, 5).Value = $r.date         # E - Date
    $ws.Cells.Item($n, 6).Value = $r.home         # F - HomeTeam
    $ws.Cells.Item($n, 7).Value = $r.away         # G - AwayTeam

    $ws.Cells.Item($n, 11).Value = $r.K           # K
    $ws.Cells.Item($n, 12).Value = $r.L           # L
    $ws.Cells.Item($n, 13).Value = $r.M           # M
    $ws.Cells.Item($n, 14).Value = $r.N           # N
    $ws.Cells.Item($n, 15).Value = $r.O           # O
    $ws.Cells.Item($n, 16).Value = $r.P           # P
    $ws.Cells.Item($n, 17).Value = $r.Q           # Q
    $ws.Cells.Item($n, 18).Value = $r.R           # R
    $ws.Cells.Item($n, 19).Value = $r.S           # S
    $ws.Cells.Item($n, 20).Value = $r.T           # T
    $ws.Cells.Item($n, 21).Value = $r.U           # U
    $ws.Cells.Item($n, 22).Value = $r.V           # V
    $ws.Cells.Item($n, 23).Value = $r.W           # W
    $ws.Cells.Item($n, 24).Value = $r.X           # X
    $ws.Cells.Item($n, 25).Value = $r.Y           # Y
    $ws.Cells.Item($n, 26).Value = $r.Z           # Z
    $ws.Cells.Item($n, 27).Value = $r.AA          # AA

    # Match formatting of the previous data rows: column A (id) uses the
    # bold/bordered/centered style, column E (Date) uses the date number format.
    $ws.Range("A145").Copy()
    $ws.Range("A$n").PasteSpecial(-4122)

    $ws.Range("E145").Copy()
    $ws.Range("E$n").PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
